# AutoCommit_17 октября 2023 г. 9:12:47_SibNout2023
# Mark several more "ок" (attendance/homework check) cells on the
# "Sibirev I. V." sheet. Each new cell's formatting is copied from an
# existing "ок" cell in the same row so it reuses the workbook's existing
# cell style (border/alignment) instead of creating a brand-new style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: add H6 = "ок" (copy style/value from G6, which already has it)
$ws.Range("G6").Copy($ws.Range("H6"))

# Row 10: add F10 = "ок" (copy style/value from G10)
$ws.Range("G10").Copy($ws.Range("F10"))

# Row 12: E12 was blank, now "ок" (copy style/value from G12)
$ws.Range("G12").Copy($ws.Range("E12"))

# Row 14: add H14 = "ок" (copy style/value from G14)
$ws.Range("G14").Copy($ws.Range("H14"))

# Row 16: E16 was blank, now "ок" (copy style/value from G16)
$ws.Range("G16").Copy($ws.Range("E16"))

# Row 22: add H22 = "ок" (copy style/value from G22)
$ws.Range("G22").Copy($ws.Range("H22"))

# Row 27: E27 was blank, now "ок"; add H27 = "ок" (copy style/value from G27)
$ws.Range("G27").Copy($ws.Range("E27"))
$ws.Range("G27").Copy($ws.Range("H27"))

# Row 28: add H28 = "ок" (copy style/value from G28)
$ws.Range("G28").Copy($ws.Range("H28"))

# Row 29: add H29 = "ок" (copy style/value from C29, the row's own "ок" cell)
$ws.Range("C29").Copy($ws.Range("H29"))
